$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.255.02"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").Value = "1.796.33"
$ws.Range("E3").Value = "  +1.70%  "

$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "'325.76"
$ws.Range("E5").Value = "  -2.91%  "

$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").Value = "'0.4475"
$ws.Range("E7").Value = "  +14.30%  "

$ws.Range("D8").Value = "'0.3715"
$ws.Range("E8").Value = "  +9.19%  "

$ws.Range("D9").Value = "'44.67"
$ws.Range("E9").Value = "  -1.29%  "

$ws.Range("D10").Value = "'1.147"
$ws.Range("E10").Value = "  +1.76%  "

$ws.Range("D11").Value = "'0.07495"
$ws.Range("E11").Value = "  +3.48%  "

$ws.Range("D12").Value = "'22.58"
$ws.Range("E12").Value = "  +0.34%  "

$ws.Range("D13").Value = "'1.001"
$ws.Range("E13").Value = "  +0.01%  "

$ws.Range("D14").Value = "'6.265"
$ws.Range("E14").Value = "  +1.59%  "

$ws.Range("E15").Value = "  +5.64%  "

$ws.Range("D16").Value = "1.793.80"
$ws.Range("E16").Value = "  +2.01%  "

$ws.Range("D17").Value = "'0.00001086"
$ws.Range("E17").Value = "  +2.25%  "

$ws.Range("D18").Value = "'0.06742"
$ws.Range("E18").Value = "  +1.81%  "

$ws.Range("D19").Value = "'80.93"
$ws.Range("E19").Value = "  +0.67%  "

$ws.Range("E20").Value = "  +0.15%  "

$ws.Range("D21").Value = "'17.49"
$ws.Range("E21").Value = "  +2.94%  "

$ws.Range("D22").Value = "'6.304"
$ws.Range("E22").Value = "  +1.03%  "

$ws.Range("D23").Value = "28.270.04"
$ws.Range("E23").Value = "  +0.61%  "

$ws.Range("D24").Value = "'11.75"
$ws.Range("E24").Value = "  +0.45%  "

$ws.Range("D25").Value = "'2.423"
$ws.Range("E25").Value = "  +1.71%  "

$ws.Range("D27").Value = "'152.00"
$ws.Range("E27").Value = "  -1.82%  "

$ws.Range("E28").Value = "  +0.59%  "

$ws.Range("D29").Value = "1.996.96"
$ws.Range("E29").Value = "  +1.83%  "

$ws.Range("D30").Value = "'132.58"
$ws.Range("E30").Value = "  +2.32%  "

$ws.Range("D31").Value = "'1.229"
$ws.Range("E31").Value = "  -4.89%  "

$ws.Range("D32").Value = "'4.016"
$ws.Range("E32").Value = "  -1.41%  "

$ws.Range("D33").Value = "'5.802"
$ws.Range("E33").Value = "  -0.48%  "

$ws.Range("D34").Value = "'0.09360"
$ws.Range("E34").Value = "  +7.03%  "

$ws.Range("D35").Value = "'0.2335"
$ws.Range("E35").Value = "  +9.98%  "

$ws.Range("D36").Value = "'12.07"
$ws.Range("E36").Value = "  -0.26%  "

$ws.Range("D37").Value = "'0.06297"
$ws.Range("E37").Value = "  +1.13%  "

$ws.Range("D38").Value = "'0.02328"
$ws.Range("E38").Value = "  +1.26%  "

$ws.Range("D39").Value = "'5.157"
$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("D40").Value = "'0.6543"
$ws.Range("E40").Value = "  +0.27%  "

$ws.Range("B41").Value = "WEMIXTOKEN"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "'1.471"
$ws.Range("E41").Value = "  -2.07%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'8.255"
$ws.Range("E42").Value = "  +4.05%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.205"
$ws.Range("E43").Value = "  +0.25%  "

$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  +0.16%  "

$ws.Range("D45").Value = "'13.99"
$ws.Range("E45").Value = "  +1.46%  "

$ws.Range("D46").Value = "'0.6061"
$ws.Range("E46").Value = "  +0.63%  "

$ws.Range("D47").Value = "'3.770"
$ws.Range("E47").Value = "  -1.61%  "

$ws.Range("D48").Value = "'129.17"
$ws.Range("E48").Value = "  +1.43%  "

$ws.Range("E49").Value = "  +0.74%  "

$ws.Range("E50").Value = "  +1.24%  "

$ws.Range("D51").Value = "'1.154"
$ws.Range("E51").Value = "  -0.63%  "
